# Extend the Excel importer test artifact ("svein harald") with a new
# "LastMetadata" sheet, a snapshot-copy of the "Metadata" sheet, appended
# as the last (6th) worksheet and made the active tab.

$wb = $excel.ActiveWorkbook
$metadata = $wb.Worksheets.Item("Metadata")

# Add the new worksheet right after the current last sheet, then rename it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "LastMetadata"

# Reuse the existing "bold label" style (s="1") used on column A of the
# Metadata sheet, without generating any new style entries.
$metadata.Range("A1:A14").Copy() | Out-Null
$ws.Range("A1:A14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 1: role ---
$ws.Cells.Item(1,1).Value = "role"
$ws.Cells.Item(1,2).Value = "information architect"

# --- Row 2: dataModelType ---
$ws.Cells.Item(2,1).Value = "dataModelType"
$ws.Cells.Item(2,2).Value = "enterprise"

# --- Row 3: schema (value differs from Metadata: "complete") ---
$ws.Cells.Item(3,1).Value = "schema"
$ws.Cells.Item(3,2).Value = "complete"

# --- Row 4: extension ---
$ws.Cells.Item(4,1).Value = "extension"
$ws.Cells.Item(4,2).Value = "addition"

# --- Row 5: prefix ---
$ws.Cells.Item(5,1).Value = "prefix"
$ws.Cells.Item(5,2).Value = "power"

# --- Row 6: namespace ---
$ws.Cells.Item(6,1).Value = "namespace"
$ws.Cells.Item(6,2).Value = "https://purl.orgl/neat/power/"

# --- Row 7: title ---
$ws.Cells.Item(7,1).Value = "title"
$ws.Cells.Item(7,2).Value = "Power to Consumer Data Model"

# --- Row 8: description (no value) ---
$ws.Cells.Item(8,1).Value = "description"

# --- Row 9: version ---
$ws.Cells.Item(9,1).Value = "version"
$ws.Cells.Item(9,2).Value = "0.1.0"

# --- Row 10: created (date/time serial, custom format) ---
$ws.Cells.Item(10,1).Value = "created"
$ws.Cells.Item(10,2).Value = 45431.300828611107
$ws.Cells.Item(10,2).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"

# --- Row 11: updated (date/time serial, custom format) ---
$ws.Cells.Item(11,1).Value = "updated"
$ws.Cells.Item(11,2).Value = 45431.300828611107
$ws.Cells.Item(11,2).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"

# --- Row 12: creator (value differs from Metadata: no ", Svein Harald") ---
$ws.Cells.Item(12,1).Value = "creator"
$ws.Cells.Item(12,2).Value = "Jon,Emma,David,Alice"

# --- Row 13: license (no value) ---
$ws.Cells.Item(13,1).Value = "license"

# --- Row 14: rights (no value) ---
$ws.Cells.Item(14,1).Value = "rights"

# Row heights match the Metadata sheet's auto-computed row height.
for ($r = 1; $r -le 14; $r++) {
  $ws.Rows.Item($r).RowHeight = 15.6
}

# Column widths, close to the Metadata sheet's bestFit column A / custom column B.
$ws.Columns.Item(1).ColumnWidth = 15.0
$ws.Columns.Item(2).ColumnWidth = 26.67

# Selection on this sheet differs from Metadata's (B19 -> A15), and this
# sheet becomes the active (selected) tab of the workbook.
$ws.Range("A15").Select() | Out-Null
$excel.ActiveWindow.ActiveSheet.Name | Out-Null
